$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected Cypher query for the ParticipantsTab row (B2): sort the collected
# sample ids with apoc.coll.sort before joining them, per Yizhen's fix for
# \omic\ in the CDS library-selection test cases.
$newQuery = @'
Match (f)<--(g:genomic_info)
WHERE g.library_selection in ['Hybrid Selection']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p, s, apoc.coll.sort(collect(distinct(samp.sample_id))) as samples
RETURN 
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samples, ','), '') as `Samples`
ORDER BY `Participant ID`LIMIT 100
'@

$ws.Range("B2").Value = $newQuery

# The longer query text now wraps across more lines, so the row grows taller.
$ws.Rows.Item(2).RowHeight = 382.5

# Leave the selection on the cell that was edited.
$ws.Range("B2").Select()
